# Apply updated "nombre_aides" (col C) and "montant_total" (col D) figures
# for the 2020-12-30 Fonds de solidarite volet 2 regional/categorie-juridique
# data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 9;   C = "316";  D = "1042987.10" },
    @{ Row = 11;  C = "532";  D = "3836269.21" },
    @{ Row = 17;  C = "707";  D = "6573657.56" },
    @{ Row = 52;  C = "798";  D = "5186204.95" },
    @{ Row = 80;  C = "455";  D = "1479752.96" },
    @{ Row = 82;  C = "1278"; D = "10205367.06" },
    @{ Row = 94;  C = "270";  D = "738250.00" },
    @{ Row = 96;  C = "656";  D = "4504105.68" },
    @{ Row = 104; C = "1705"; D = "9801214.95" },
    @{ Row = 106; C = "1666"; D = "9105315.02" }
)

foreach ($u in $updates) {
    $cCell = $ws.Range("C" + $u.Row)
    $cCell.NumberFormat = "@"
    $cCell.Value = $u.C

    $dCell = $ws.Range("D" + $u.Row)
    $dCell.NumberFormat = "@"
    $dCell.Value = $u.D
}

Write-Output "Updated $($updates.Count) rows"
